# Update "想去人数" (F column) figures across the 展览, 演出, and 全部类型 sheets
# to match the freshly generated gh-pages output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 75
$ws1.Range("F6").Value  = 33
$ws1.Range("F7").Value  = 2626
$ws1.Range("F8").Value  = 1141
$ws1.Range("F9").Value  = 229
$ws1.Range("F10").Value = 90
$ws1.Range("F11").Value = 5892
$ws1.Range("F13").Value = 230
$ws1.Range("F14").Value = 577
$ws1.Range("F15").Value = 11568
$ws1.Range("F16").Value = 11745
$ws1.Range("F18").Value = 74
$ws1.Range("F21").Value = 60
$ws1.Range("F22").Value = 29

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 6

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 75
$ws4.Range("F6").Value  = 33
$ws4.Range("F7").Value  = 2626
$ws4.Range("F8").Value  = 6
$ws4.Range("F9").Value  = 1141
$ws4.Range("F10").Value = 229
$ws4.Range("F11").Value = 90
$ws4.Range("F12").Value = 5892
$ws4.Range("F14").Value = 230
$ws4.Range("F15").Value = 577
$ws4.Range("F16").Value = 11568
$ws4.Range("F17").Value = 11745
$ws4.Range("F19").Value = 74
$ws4.Range("F22").Value = 60
$ws4.Range("F23").Value = 29
